# Quarterly update: drop oldest period column, shift remaining periods left,
# and append the new "3 ماهه منتهی به 1401/11" / "1401-12-29" period at column M.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# --- Row 8: financial period headers ---
$ws.Range("D8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("E8").Value = "3 ماهه منتهی به 1399/08"
$ws.Range("F8").Value = "6 ماهه منتهی به 1400/06"
$ws.Range("G8").Value = "9 ماهه منتهی به 1400/09"
$ws.Range("H8").Value = "12 ماهه منتهی به 1400/08"
$ws.Range("I8").Value = "3 ماهه منتهی به 1400/11"
$ws.Range("J8").Value = "6 ماهه منتهی به 1401/02"
$ws.Range("K8").Value = "9 ماهه منتهی به 1401/05"
$ws.Range("L8").Value = "12 ماهه منتهی به 1401/08"
$ws.Range("M8").Value = "3 ماهه منتهی به 1401/11"

# --- Row 9: publish dates ---
$ws.Range("D9").Value = "1400-10-08 (15)"
$ws.Range("E9").Value = "1399-09-30"
$ws.Range("F9").Value = "1400-07-29"
$ws.Range("G9").ClearContents()
$ws.Range("H9").Value = "1401-10-28 (10)"
$ws.Range("I9").Value = "1400-12-28"
$ws.Range("J9").Value = "1401-04-29 (2)"
$ws.Range("K9").Value = "1401-06-30"
$ws.Range("L9").Value = "1401-10-28 (3)"
$ws.Range("M9").Value = "1401-12-29"

# Row 11: فروش
$ws.Range("D11").Value = 1994706
$ws.Range("E11").Value = 761556
$ws.Range("F11").Value = 2037897
$ws.Range("G11").Value = "-"
$ws.Range("H11").Value = 3049033
$ws.Range("I11").Value = 1803801
$ws.Range("J11").Value = 3883709
$ws.Range("K11").Value = 5718026
$ws.Range("L11").Value = 8001536
$ws.Range("M11").Value = 2241847

# Row 12: بهای تمام شده کالای فروش رفته
$ws.Range("D12").Value = -968950
$ws.Range("E12").Value = -339263
$ws.Range("F12").Value = -1177381
$ws.Range("G12").Value = "-"
$ws.Range("H12").Value = -1724946
$ws.Range("I12").Value = -567700
$ws.Range("J12").Value = -1328999
$ws.Range("K12").Value = -2127202
$ws.Range("L12").Value = -2980516
$ws.Range("M12").Value = -840364

# Row 13: سود (زیان) ناخالص
$ws.Range("D13").Value = 1025756
$ws.Range("E13").Value = 422293
$ws.Range("F13").Value = 860516
$ws.Range("G13").Value = "-"
$ws.Range("H13").Value = 1324087
$ws.Range("I13").Value = 1236101
$ws.Range("J13").Value = 2554710
$ws.Range("K13").Value = 3590824
$ws.Range("L13").Value = 5021020
$ws.Range("M13").Value = 1401483

# Row 14: هزینه های عمومی, اداری و تشکیلاتی
$ws.Range("D14").Value = -96381
$ws.Range("E14").Value = -31248
$ws.Range("F14").Value = -151528
$ws.Range("G14").Value = "-"
$ws.Range("H14").Value = -253507
$ws.Range("I14").Value = -95724
$ws.Range("J14").Value = -238494
$ws.Range("K14").Value = -430022
$ws.Range("L14").Value = -496855
$ws.Range("M14").Value = -151717

# Row 15: هزینه کاهش ارزش دریافتنی‌‏ها (هزینه استثنایی)
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = "-"
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 0

# Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی
$ws.Range("D16").Value = 83739
$ws.Range("E16").Value = 48574
$ws.Range("F16").Value = 70209
$ws.Range("G16").Value = "-"
$ws.Range("H16").Value = 119005
$ws.Range("I16").Value = 40771
$ws.Range("J16").Value = 52335
$ws.Range("K16").Value = 127529
$ws.Range("L16").Value = 189764
$ws.Range("M16").Value = 163473

# Row 17: سود (زیان) عملیاتی
$ws.Range("D17").Value = 1013114
$ws.Range("E17").Value = 439619
$ws.Range("F17").Value = 779197
$ws.Range("G17").Value = "-"
$ws.Range("H17").Value = 1189585
$ws.Range("I17").Value = 1181148
$ws.Range("J17").Value = 2368551
$ws.Range("K17").Value = 3288331
$ws.Range("L17").Value = 4713929
$ws.Range("M17").Value = 1413239

# Row 18: هزینه های مالی
$ws.Range("D18").Value = -32821
$ws.Range("E18").Value = -14780
$ws.Range("F18").Value = -39598
$ws.Range("G18").Value = "-"
$ws.Range("H18").Value = -56214
$ws.Range("I18").Value = -13975
$ws.Range("J18").Value = -13028
$ws.Range("K18").Value = -13028
$ws.Range("L18").Value = -25407
$ws.Range("M18").Value = 0

# Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی
$ws.Range("D19").Value = 50013
$ws.Range("E19").Value = 24755
$ws.Range("F19").Value = 40831
$ws.Range("G19").Value = "-"
$ws.Range("H19").Value = 42428
$ws.Range("I19").Value = 8499
$ws.Range("J19").Value = 26499
$ws.Range("K19").Value = 31725
$ws.Range("L19").Value = 113030
$ws.Range("M19").Value = 36211

# Row 20: سود (زیان) خالص عملیات در حال تداوم قبل از مالیات
$ws.Range("D20").Value = 1030306
$ws.Range("E20").Value = 449594
$ws.Range("F20").Value = 780430
$ws.Range("G20").Value = "-"
$ws.Range("H20").Value = 1175799
$ws.Range("I20").Value = 1175672
$ws.Range("J20").Value = 2382022
$ws.Range("K20").Value = 3307028
$ws.Range("L20").Value = 4801552
$ws.Range("M20").Value = 1449450

# Row 21: مالیات
$ws.Range("D21").Value = -106773
$ws.Range("E21").Value = -28939
$ws.Range("F21").Value = -68807
$ws.Range("G21").Value = "-"
$ws.Range("H21").Value = -143742
$ws.Range("I21").Value = -158181
$ws.Range("J21").Value = -354869
$ws.Range("K21").Value = -479943
$ws.Range("L21").Value = -607778
$ws.Range("M21").Value = -226285

# Row 22: سود (زیان) خالص عملیات در حال تداوم
$ws.Range("D22").Value = 923533
$ws.Range("E22").Value = 420655
$ws.Range("F22").Value = 711623
$ws.Range("G22").Value = "-"
$ws.Range("H22").Value = 1032057
$ws.Range("I22").Value = 1017491
$ws.Range("J22").Value = 2027153
$ws.Range("K22").Value = 2827085
$ws.Range("L22").Value = 4193774
$ws.Range("M22").Value = 1223165

# Row 23: سود (زیان) عملیات متوقف شده پس از اثر مالیاتی
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = "-"
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 0

# Row 24: سود (زیان) خالص
$ws.Range("D24").Value = 923533
$ws.Range("E24").Value = 420655
$ws.Range("F24").Value = 711623
$ws.Range("G24").Value = "-"
$ws.Range("H24").Value = 1032057
$ws.Range("I24").Value = 1017491
$ws.Range("J24").Value = 2027153
$ws.Range("K24").Value = 2827085
$ws.Range("L24").Value = 4193774
$ws.Range("M24").Value = 1223165

# Row 25: سود هر سهم پس از کسر مالیات
$ws.Range("D25").Value = 1847
$ws.Range("E25").Value = 841
$ws.Range("F25").Value = 1423
$ws.Range("G25").Value = "-"
$ws.Range("H25").Value = 1032
$ws.Range("I25").Value = 2035
$ws.Range("J25").Value = 4054
$ws.Range("K25").Value = 2827
$ws.Range("L25").Value = 4194
$ws.Range("M25").Value = 1223

# Row 26: سرمایه
$ws.Range("D26").Value = 500000
$ws.Range("E26").Value = 500000
$ws.Range("F26").Value = 500000
$ws.Range("G26").Value = "-"
$ws.Range("H26").Value = 1000000
$ws.Range("I26").Value = 500000
$ws.Range("J26").Value = 500000
$ws.Range("K26").Value = 1000000
$ws.Range("L26").Value = 1000000
$ws.Range("M26").Value = 1000000

# Row 27: سود هر سهم بر اساس آخرین سرمایه
$ws.Range("D27").Value = 924
$ws.Range("E27").Value = 421
$ws.Range("F27").Value = 712
$ws.Range("G27").Value = "-"
$ws.Range("H27").Value = 1032
$ws.Range("I27").Value = 1017
$ws.Range("J27").Value = 2027
$ws.Range("K27").Value = 2827
$ws.Range("L27").Value = 4194
$ws.Range("M27").Value = 1223
